$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 132: Fast-forwarding Flora
$ws.Cells.Item(132, 8).Value = 12828231
$ws.Cells.Item(132, 9).Value = 12828231
$ws.Cells.Item(132, 11).Value = 38484693
$ws.Cells.Item(132, 13).Value = -38482163

# Row 137: Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 2670.8708
$ws.Cells.Item(137, 9).Value = 2587.775
$ws.Cells.Item(137, 11).Value = 7763.325000000001
$ws.Cells.Item(137, 13).Value = -5213.325000000001

# Row 138: All-night Crafting
$ws.Cells.Item(138, 8).Value = 4017.081
$ws.Cells.Item(138, 9).Value = 1521.8572
$ws.Cells.Item(138, 10).Value = 6256.385
$ws.Cells.Item(138, 11).Value = 4565.571599999999
$ws.Cells.Item(138, 12).Value = 18769.155
$ws.Cells.Item(138, 13).Value = 574.4284000000007
$ws.Cells.Item(138, 14).Value = -29049.155

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Cells.Item(32, 8).Value = 4699508.5
$ws.Cells.Item(32, 9).Value = 4477.5
$ws.Cells.Item(32, 10).Value = 19613136
$ws.Cells.Item(32, 11).Value = 4477.5
$ws.Cells.Item(32, 12).Value = 19613136
$ws.Cells.Item(32, 13).Value = -4190.5
$ws.Cells.Item(32, 14).Value = -19613710

# Row 110: Scheduled Maintenance
$ws.Cells.Item(110, 8).Value = 2545.5557
$ws.Cells.Item(110, 9).Value = 1753
$ws.Cells.Item(110, 10).Value = 5319.5
$ws.Cells.Item(110, 11).Value = 1753
$ws.Cells.Item(110, 12).Value = 5319.5
$ws.Cells.Item(110, 13).Value = 292
$ws.Cells.Item(110, 14).Value = -9409.5

# Row 122: Haste for High Durium
$ws.Cells.Item(122, 8).Value = 1438.8
$ws.Cells.Item(122, 9).Value = 1058.4
$ws.Cells.Item(122, 10).Value = 2580
$ws.Cells.Item(122, 11).Value = 3175.2
$ws.Cells.Item(122, 12).Value = 7740
$ws.Cells.Item(122, 13).Value = -725.2000000000003
$ws.Cells.Item(122, 14).Value = -12640

# Row 132: Don't Bore Me, Ore Me
$ws.Cells.Item(132, 8).Value = 879764.5600000001
$ws.Cells.Item(132, 9).Value = 1683.5834
$ws.Cells.Item(132, 10).Value = 3098074.5
$ws.Cells.Item(132, 11).Value = 5050.7502
$ws.Cells.Item(132, 12).Value = 9294223.5
$ws.Cells.Item(132, 13).Value = -2520.7502
$ws.Cells.Item(132, 14).Value = -9299283.5

$ws = $wb.Worksheets.Item("BSM")
# Row 107: The Gold Experience
$ws.Cells.Item(107, 8).Value = 16668182
$ws.Cells.Item(107, 9).Value = 26316792
$ws.Cells.Item(107, 10).Value = 2401.0908
$ws.Cells.Item(107, 11).Value = 26316792
$ws.Cells.Item(107, 12).Value = 2401.0908
$ws.Cells.Item(107, 13).Value = -26314872
$ws.Cells.Item(107, 14).Value = -6241.0908

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting
$ws.Cells.Item(58, 8).Value = 31250888
$ws.Cells.Item(58, 9).Value = 38462216
$ws.Cells.Item(58, 10).Value = 1799.8334
$ws.Cells.Item(58, 11).Value = 38462216
$ws.Cells.Item(58, 12).Value = 1799.8334
$ws.Cells.Item(58, 13).Value = -38462013
$ws.Cells.Item(58, 14).Value = -2205.8334

# Row 122: Timber of Tenkonto
$ws.Cells.Item(122, 8).Value = 17859232
$ws.Cells.Item(122, 9).Value = 22729384
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 68188152
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -68185702
$ws.Cells.Item(122, 14).Value = -10900

# Row 136: Turali Quality
$ws.Cells.Item(136, 8).Value = 31250888
$ws.Cells.Item(136, 9).Value = 38462216
$ws.Cells.Item(136, 10).Value = 1799.8334
$ws.Cells.Item(136, 11).Value = 115386648
$ws.Cells.Item(136, 12).Value = 5399.5002
$ws.Cells.Item(136, 13).Value = -115384098
$ws.Cells.Item(136, 14).Value = -10499.5002

$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face
$ws.Cells.Item(68, 8).Value = 2858.2593
$ws.Cells.Item(68, 9).Value = 812.0303
$ws.Cells.Item(68, 10).Value = 6073.7617
$ws.Cells.Item(68, 11).Value = 2436.0909
$ws.Cells.Item(68, 12).Value = 18221.2851
$ws.Cells.Item(68, 13).Value = -1625.0909
$ws.Cells.Item(68, 14).Value = -19843.2851

# Row 71: No Margarine of Error (L)
$ws.Cells.Item(71, 8).Value = 2858.2593
$ws.Cells.Item(71, 9).Value = 812.0303
$ws.Cells.Item(71, 10).Value = 6073.7617
$ws.Cells.Item(71, 11).Value = 7308.2727
$ws.Cells.Item(71, 12).Value = 54663.8553
$ws.Cells.Item(71, 13).Value = -3252.2727
$ws.Cells.Item(71, 14).Value = -62775.8553

# Row 87: Soup That Eats Like a Knight
$ws.Cells.Item(87, 8).Value = 1906
$ws.Cells.Item(87, 9).Value = 1906
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 5718
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Value = -4470
$ws.Cells.Item(87, 14).ClearContents()

# Row 90: Like Ma Used to Make (L)
$ws.Cells.Item(90, 8).Value = 1906
$ws.Cells.Item(90, 9).Value = 1906
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 11).Value = 17154
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 13).Value = -10914
$ws.Cells.Item(90, 14).ClearContents()

# Row 131: The Mountain Steeped
$ws.Cells.Item(131, 8).Value = 840.15
$ws.Cells.Item(131, 9).Value = 426.1875
$ws.Cells.Item(131, 10).Value = 919
$ws.Cells.Item(131, 11).Value = 1278.5625
$ws.Cells.Item(131, 12).Value = 2757
$ws.Cells.Item(131, 13).Value = 3761.4375
$ws.Cells.Item(131, 14).Value = -12837

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence
$ws.Cells.Item(122, 8).Value = 100042560
$ws.Cells.Item(122, 9).Value = 125050700
$ws.Cells.Item(122, 10).Value = 10000
$ws.Cells.Item(122, 11).Value = 375152100
$ws.Cells.Item(122, 12).Value = 30000
$ws.Cells.Item(122, 13).Value = -375149650
$ws.Cells.Item(122, 14).Value = -34900

# Row 132: On Board for Lar
$ws.Cells.Item(132, 8).Value = 8254.579
$ws.Cells.Item(132, 9).Value = 2233.9
$ws.Cells.Item(132, 10).Value = 14944.223
$ws.Cells.Item(132, 11).Value = 6701.700000000001
$ws.Cells.Item(132, 12).Value = 44832.669
$ws.Cells.Item(132, 13).Value = -4171.700000000001
$ws.Cells.Item(132, 14).Value = -49892.669

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad
$ws.Cells.Item(40, 8).Value = 83336184
$ws.Cells.Item(40, 9).Value = 3000
$ws.Cells.Item(40, 11).Value = 3000
$ws.Cells.Item(40, 13).Value = -2864

# Row 100: Tiger in the Sack
$ws.Cells.Item(100, 8).Value = 3576.8462
$ws.Cells.Item(100, 9).Value = 3100
$ws.Cells.Item(100, 10).Value = 3788.7778
$ws.Cells.Item(100, 11).Value = 3100
$ws.Cells.Item(100, 12).Value = 3788.7778
$ws.Cells.Item(100, 13).Value = -2559
$ws.Cells.Item(100, 14).Value = -4870.7778

# Row 122: Hell on Leather
$ws.Cells.Item(122, 8).Value = 8861.308000000001
$ws.Cells.Item(122, 9).Value = 12333.556
$ws.Cells.Item(122, 10).Value = 1048.75
$ws.Cells.Item(122, 11).Value = 37000.66800000001
$ws.Cells.Item(122, 12).Value = 3146.25
$ws.Cells.Item(122, 13).Value = -34550.66800000001
$ws.Cells.Item(122, 14).Value = -8046.25

# Row 132: Tenets of Tanning
$ws.Cells.Item(132, 8).Value = 27787432
$ws.Cells.Item(132, 9).Value = 52634740
$ws.Cells.Item(132, 10).Value = 16912.234
$ws.Cells.Item(132, 11).Value = 157904220
$ws.Cells.Item(132, 12).Value = 50736.702
$ws.Cells.Item(132, 13).Value = -157901690
$ws.Cells.Item(132, 14).Value = -55796.702

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches
$ws.Cells.Item(81, 8).Value = 1195.1111
$ws.Cells.Item(81, 9).Value = 1030.4
$ws.Cells.Item(81, 10).Value = 1401
$ws.Cells.Item(81, 11).Value = 2060.8
$ws.Cells.Item(81, 12).Value = 2802
$ws.Cells.Item(81, 13).Value = -999.8000000000002
$ws.Cells.Item(81, 14).Value = -4924

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Cells.Item(84, 8).Value = 1195.1111
$ws.Cells.Item(84, 9).Value = 1030.4
$ws.Cells.Item(84, 10).Value = 1401
$ws.Cells.Item(84, 11).Value = 10304
$ws.Cells.Item(84, 12).Value = 14010
$ws.Cells.Item(84, 13).Value = -5000
$ws.Cells.Item(84, 14).Value = -24618

# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 46607.54
$ws.Cells.Item(132, 9).Value = 67237
$ws.Cells.Item(132, 10).Value = 13600.4
$ws.Cells.Item(132, 11).Value = 201711
$ws.Cells.Item(132, 12).Value = 40801.2
$ws.Cells.Item(132, 13).Value = -199181
$ws.Cells.Item(132, 14).Value = -45861.2

# Row 136: Weaving the Envelope
$ws.Cells.Item(136, 8).Value = 15627779
$ws.Cells.Item(136, 9).Value = 45456830
$ws.Cells.Item(136, 10).Value = 3037.1428
$ws.Cells.Item(136, 11).Value = 136370490
$ws.Cells.Item(136, 12).Value = 9111.428400000001
$ws.Cells.Item(136, 13).Value = -136367940
$ws.Cells.Item(136, 14).Value = -14211.4284
